$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to already-populated text cells (simple value replacement) ---
$ws.Range("B4").Value  = "Supervised Internship"
$ws.Range("C4").Value  = "Supervised Internship"

$ws.Range("B7").Value  = "195 h   (    Estágio: 195 h         )"
$ws.Range("C7").Value  = "195 h   (    Estágio: 195 h         )"

$ws.Range("B10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C10").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
$ws.Range("C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"

# --- Updates to cells whose new value looks like a number/date: force text
#     entry (leading apostrophe) so Excel keeps them as plain text instead of
#     auto-converting to a numeric / date value, matching the source data. ---
$ws.Range("B6").Value  = "'6"
$ws.Range("C6").Value  = "'6"

$ws.Range("B8").Value  = "'01/01/2023"
$ws.Range("C8").Value  = "'01/01/2023"

$ws.Range("B15").Value = "'01/01/2023"
$ws.Range("C15").Value = "'01/01/2023"

# --- Newly populated cells (previously empty): copy the number format from a
#     neighboring populated cell in the same column first, so the new cell
#     picks up the same wrap/alignment style instead of the plain column
#     default, then set its value. ---
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B11").Value = "Offer the opportunity to carry out professional training in a company or research institution, under the supervision of a professor from the Materials Engineering Department at EEL. Complement the general curricular training and psychologically and socially adapt the student to his/her future professional activity."

$ws.Range("C10").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = "Offer the opportunity to carry out professional training in a company or research institution, under the supervision of a professor from the Materials Engineering Department at EEL. Complement the general curricular training and psychologically and socially adapt the student to his/her future professional activity."

$ws.Range("B13").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("B14").Value = "Participation in the selection process or indication of an institution to carry out an internship. Submission of the specific work plan. Conducting the internship and delivering the internship report."

$ws.Range("C13").Copy() | Out-Null
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = "Participation in the selection process or indication of an institution to carry out an internship. Submission of the specific work plan. Conducting the internship and delivering the internship report."

$ws.Range("B18").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = "Student participation in the selection process of companies, research institutions or in the academic sector. The internship will be carried out under the supervision of a professor appointed by the Physical Engineering Course Committee. The content will be established in the Work Plan between the supervisor responsible for the Internship and the supervising professor. Presentation of a final report on the activities carried out in the internship."

$ws.Range("C18").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "Student participation in the selection process of companies, research institutions or in the academic sector. The internship will be carried out under the supervision of a professor appointed by the Physical Engineering Course Committee. The content will be established in the Work Plan between the supervisor responsible for the Internship and the supervising professor. Presentation of a final report on the activities carried out in the internship."
